$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append missing Maxwell et al. 2023 studies (study_id in column A, doi in column B)
# Row 23
$ws.Cells.Item(23, 1).Value = "Yando_et_al_2016"
$ws.Cells.Item(23, 2).Value = "10.1111/1365-2745.12571"

# Row 24
$ws.Cells.Item(24, 1).Value = "Gao_et_al_2016"
$ws.Cells.Item(24, 2).Value = "10.1016/j.ecoleng.2016.06.088"

# Row 25 (doi entered before study_id)
$ws.Cells.Item(25, 2).Value = "10.1016/j.ecoleng.2017.05.041"
$ws.Cells.Item(25, 1).Value = "Liu_et_al_2017"

# Row 26
$ws.Cells.Item(26, 1).Value = "Fu_et_al_2021"
$ws.Cells.Item(26, 2).Value = "10.1111/gcb.15348"

$ws.Range("A26").Select()
